$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Activate()

# Updated fitting parameters (r_s_star, h_p_star)
$ws.Range("J2").Value = 0.01525699
$ws.Range("K2").Value = 0.1115

# Leave the selection where the user ended up after editing K2
$ws.Range("K3").Select()
